$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.460.13'
$ws.Range("E2").Value = '  -1.08%  '

$ws.Range("D3").Value = '3.100.65'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.31'
$ws.Range("E5").Value = '  +6.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '624.31'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  +5.88%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.350'
$ws.Range("E8").Value = '  -7.22%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").Value = '3.381.77'
$ws.Range("E10").Value = '  +8.16%  '

$ws.Range("E11").Value = '  -4.32%  '

$ws.Range("E12").Value = '  +3.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.93'
$ws.Range("E13").Value = '  +2.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000241'
$ws.Range("E14").Value = '  -4.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.55'
$ws.Range("E15").Value = '  +2.04%  '

$ws.Range("D16").Value = '90.277.90'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = '3.684.04'
$ws.Range("E17").Value = '  -0.69%  '

$ws.Range("D18").Value = '3.065.39'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.66'
$ws.Range("E19").Value = '  -2.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.19'
$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000210'
$ws.Range("E21").Value = '  -7.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.07'
$ws.Range("E22").Value = '  +2.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("E23").Value = '  +1.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.35'
$ws.Range("E24").Value = '  +3.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.07'
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.69'
$ws.Range("E26").Value = '  +4.39%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.16'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.21'
$ws.Range("E30").Value = '  +2.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.158'
$ws.Range("E31").Value = '  -7.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.19'
$ws.Range("E32").Value = '  +14.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.198'
$ws.Range("E33").Value = '  +27.93%  '

$ws.Range("E34").Value = '  +3.66%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '504.97'
$ws.Range("E35").Value = '  -4.88%  '

$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.73'
$ws.Range("E36").Value = '  -3.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("E37").Value = '  +2.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.97'

$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.419'
$ws.Range("E40").Value = '  +8.91%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.17'

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0849'
$ws.Range("E43").Value = '  +7.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.737'
$ws.Range("E44").Value = '  -17.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.26'
$ws.Range("E45").Value = '  +34.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.95'
$ws.Range("E46").Value = '  +0.79%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.94'
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.691'
$ws.Range("E48").Value = '  +10.40%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.84'
$ws.Range("E49").Value = '  +1.64%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.50'
$ws.Range("E50").Value = '  +7.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.34'
$ws.Range("E51").Value = '  +3.01%  '
